# Updates the cryptos price/volume snapshot (GitHub Actions scheduled refresh).
# Values are written through a TEXT-FORMULA + copy/paste-values round trip so that
# numeric-looking strings (e.g. "49.10", "0.06900") stay literal text instead of
# being auto-coerced to numbers (which would drop trailing zeros / use sci notation)
# and so that no cell style/number-format gets attached as a side effect.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: "27.963.59" -> "27.978.16"
$ws.Range("D2").Formula = "=""27.978.16"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

# E2: "  -2.82%  " -> "  -2.60%  "
$ws.Range("E2").Formula = "=""  -2.60%  """
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

# D3: "1.888.59" -> "1.890.04"
$ws.Range("D3").Formula = "=""1.890.04"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

# E3: "  -3.64%  " -> "  -3.34%  "
$ws.Range("E3").Formula = "=""  -3.34%  """
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

# E4: "  -0.89%  " -> "  -1.27%  "
$ws.Range("E4").Formula = "=""  -1.27%  """
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

# D5: "325.86" -> "326.56"
$ws.Range("D5").Formula = "=""326.56"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)

# E5: "  +1.22%  " -> "  +1.46%  "
$ws.Range("E5").Formula = "=""  +1.46%  """
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)

# E6: "  -0.87%  " -> "  -1.07%  "
$ws.Range("E6").Formula = "=""  -1.07%  """
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)

# D7: "0.4571" -> "0.4589"
$ws.Range("D7").Formula = "=""0.4589"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)

# E7: "  -4.10%  " -> "  -3.64%  "
$ws.Range("E7").Formula = "=""  -3.64%  """
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)

# D8: "0.3935" -> "0.3938"
$ws.Range("D8").Formula = "=""0.3938"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)

# E8: "  -2.40%  " -> "  -2.12%  "
$ws.Range("E8").Formula = "=""  -2.12%  """
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)

# D9: "50.40" -> "49.10"
$ws.Range("D9").Formula = "=""49.10"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)

# E9: "  -6.59%  " -> "  -9.30%  "
$ws.Range("E9").Formula = "=""  -9.30%  """
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)

# D10: "0.08201" -> "0.08239"
$ws.Range("D10").Formula = "=""0.08239"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)

# E10: "  -3.28%  " -> "  -2.62%  "
$ws.Range("E10").Formula = "=""  -2.62%  """
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)

# D11: "1.036" -> "1.038"
$ws.Range("D11").Formula = "=""1.038"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)

# E11: "  -2.37%  " -> "  -1.99%  "
$ws.Range("E11").Formula = "=""  -1.99%  """
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)

# D12: "21.78" -> "21.79"
$ws.Range("D12").Formula = "=""21.79"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)

# E12: "  -2.70%  " -> "  -2.22%  "
$ws.Range("E12").Formula = "=""  -2.22%  """
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

# D13: "1.911.52" -> "1.893.03"
$ws.Range("D13").Formula = "=""1.893.03"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)

# E13: "  -4.03%  " -> "  -3.51%  "
$ws.Range("E13").Formula = "=""  -3.51%  """
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)

# D14: "7.331" -> "7.334"
$ws.Range("D14").Formula = "=""7.334"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)

# E14: "  -3.82%  " -> "  -3.46%  "
$ws.Range("E14").Formula = "=""  -3.46%  """
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)

# D15: "5.977" -> "5.983"
$ws.Range("D15").Formula = "=""5.983"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)

# E15: "  -3.99%  " -> "  -3.75%  "
$ws.Range("E15").Formula = "=""  -3.75%  """
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)

# D16: "1.004" -> "1.002"
$ws.Range("D16").Formula = "=""1.002"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)

# E16: "  -1.01%  " -> "  -1.31%  "
$ws.Range("E16").Formula = "=""  -1.31%  """
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)

# D17: "89.39" -> "89.28"
$ws.Range("D17").Formula = "=""89.28"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)

# E17: "  -0.01%  " -> "  +0.25%  "
$ws.Range("E17").Formula = "=""  +0.25%  """
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)

# E18: "  -1.66%  " -> "  -1.78%  "
$ws.Range("E18").Formula = "=""  -1.78%  """
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)

# D19: "0.06585" -> "0.06572"
$ws.Range("D19").Formula = "=""0.06572"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)

# E19: "  -0.43%  " -> "  -0.93%  "
$ws.Range("E19").Formula = "=""  -0.93%  """
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)

# D20: "17.48" -> "17.52"
$ws.Range("D20").Formula = "=""17.52"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)

# E20: "  -6.03%  " -> "  -5.86%  "
$ws.Range("E20").Formula = "=""  -5.86%  """
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)

# D21: "1.002" -> "1.001"
$ws.Range("D21").Formula = "=""1.001"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)

# E21: "  -0.95%  " -> "  -1.04%  "
$ws.Range("E21").Formula = "=""  -1.04%  """
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)

# D22: "5.637" -> "5.636"
$ws.Range("D22").Formula = "=""5.636"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)

# E22: "  -2.66%  " -> "  -2.57%  "
$ws.Range("E22").Formula = "=""  -2.57%  """
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)

# D23: "27.942.42" -> "27.965.29"
$ws.Range("D23").Formula = "=""27.965.29"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)

# E23: "  -2.98%  " -> "  -2.74%  "
$ws.Range("E23").Formula = "=""  -2.74%  """
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)

# D24: "11.06" -> "11.08"
$ws.Range("D24").Formula = "=""11.08"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)

# E24: "  -4.10%  " -> "  -3.80%  "
$ws.Range("E24").Formula = "=""  -3.80%  """
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)

# D25: "2.302" -> "2.304"
$ws.Range("D25").Formula = "=""2.304"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)

# E25: "  +0.26%  " -> "  +0.29%  "
$ws.Range("E25").Formula = "=""  +0.29%  """
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)

# D26: "2.101.56" -> "2.131.21"
$ws.Range("D26").Formula = "=""2.131.21"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)

# E26: "  -5.47%  " -> "  -3.74%  "
$ws.Range("E26").Formula = "=""  -3.74%  """
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)

# D27: "154.14" -> "154.26"
$ws.Range("D27").Formula = "=""154.26"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)

# E27: "  -0.16%  " -> "  -0.10%  "
$ws.Range("E27").Formula = "=""  -0.10%  """
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)

# D28: "19.88" -> "19.91"
$ws.Range("D28").Formula = "=""19.91"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)

# E28: "  -1.50%  " -> "  -1.24%  "
$ws.Range("E28").Formula = "=""  -1.24%  """
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)

# D29: "5.702" -> "5.710"
$ws.Range("D29").Formula = "=""5.710"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)

# E29: "  -3.72%  " -> "  -3.39%  "
$ws.Range("E29").Formula = "=""  -3.39%  """
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)

# D30: "2.106" -> "2.108"
$ws.Range("D30").Formula = "=""2.108"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)

# E30: "  -2.32%  " -> "  -1.96%  "
$ws.Range("E30").Formula = "=""  -1.96%  """
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)

# D31: "123.87" -> "123.60"
$ws.Range("D31").Formula = "=""123.60"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)

# E31: "  -0.06%  " -> "  +0.03%  "
$ws.Range("E31").Formula = "=""  +0.03%  """
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)

# D32: "0.09527" -> "0.09546"
$ws.Range("D32").Formula = "=""0.09546"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)

# E32: "  -0.52%  " -> "  -0.20%  "
$ws.Range("E32").Formula = "=""  -0.20%  """
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)

# D33: "0.9588" -> "0.9591"
$ws.Range("D33").Formula = "=""0.9591"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)

# E33: "  -4.50%  " -> "  -4.20%  "
$ws.Range("E33").Formula = "=""  -4.20%  """
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)

# D34: "1.476" -> "1.475"
$ws.Range("D34").Formula = "=""1.475"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)

# E34: "  +1.89%  " -> "  +2.78%  "
$ws.Range("E34").Formula = "=""  +2.78%  """
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)

# D35: "3.622" -> "3.630"
$ws.Range("D35").Formula = "=""3.630"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)

# E35: "  -1.14%  " -> "  -0.95%  "
$ws.Range("E35").Formula = "=""  -0.95%  """
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)

# D36: "5.462" -> "5.470"
$ws.Range("D36").Formula = "=""5.470"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)

# E36: "  -3.75%  " -> "  -3.34%  "
$ws.Range("E36").Formula = "=""  -3.34%  """
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)

# E37: "  -0.76%  " -> "  -0.74%  "
$ws.Range("E37").Formula = "=""  -0.74%  """
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)

# D38: "0.02276" -> "0.02280"
$ws.Range("D38").Formula = "=""0.02280"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)

# E38: "  -2.76%  " -> "  -2.79%  "
$ws.Range("E38").Formula = "=""  -2.79%  """
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)

# D39: "0.06100" -> "0.06110"
$ws.Range("D39").Formula = "=""0.06110"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)

# E39: "  -1.87%  " -> "  -1.62%  "
$ws.Range("E39").Formula = "=""  -1.62%  """
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)

# D40: "8.601" -> "8.595"
$ws.Range("D40").Formula = "=""8.595"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)

# E40: "  -1.57%  " -> "  -1.45%  "
$ws.Range("E40").Formula = "=""  -1.45%  """
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)

# D41: "0.6103" -> "0.6109"
$ws.Range("D41").Formula = "=""0.6109"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)

# E41: "  -1.63%  " -> "  -1.39%  "
$ws.Range("E41").Formula = "=""  -1.39%  """
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)

# E42: "  -0.83%  " -> "  -0.95%  "
$ws.Range("E42").Formula = "=""  -0.95%  """
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)

# D43: "10.73" -> "10.75"
$ws.Range("D43").Formula = "=""10.75"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)

# E43: "  -3.26%  " -> "  -2.72%  "
$ws.Range("E43").Formula = "=""  -2.72%  """
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)

# D44: "0.1895" -> "0.1898"
$ws.Range("D44").Formula = "=""0.1898"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)

# E44: "  -1.16%  " -> "  -0.88%  "
$ws.Range("E44").Formula = "=""  -0.88%  """
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)

# D45: "1.307" -> "1.313"
$ws.Range("D45").Formula = "=""1.313"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)

# E45: "  -1.55%  " -> "  -1.32%  "
$ws.Range("E45").Formula = "=""  -1.32%  """
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)

# B46: "Decentraland" -> "EnergySwap"
$ws.Range("B46").Formula = "=""EnergySwap"""
$ws.Range("B46").Copy()
$ws.Range("B46").PasteSpecial(-4163)

# C46: "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana" -> "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C46").Formula = "=""https://coinranking.com/coin/SbWqqTui-+energyswap-ens"""
$ws.Range("C46").Copy()
$ws.Range("C46").PasteSpecial(-4163)

# D46: "0.5813" -> "12.75"
$ws.Range("D46").Formula = "=""12.75"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)

# E46: "  -2.17%  " -> "  -0.90%  "
$ws.Range("E46").Formula = "=""  -0.90%  """
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)

# B47: "EnergySwap" -> "Decentraland"
$ws.Range("B47").Formula = "=""Decentraland"""
$ws.Range("B47").Copy()
$ws.Range("B47").PasteSpecial(-4163)

# C47: "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" -> "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C47").Formula = "=""https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"""
$ws.Range("C47").Copy()
$ws.Range("C47").PasteSpecial(-4163)

# D47: "12.75" -> "0.5807"
$ws.Range("D47").Formula = "=""0.5807"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)

# E47: "  -1.32%  " -> "  -1.99%  "
$ws.Range("E47").Formula = "=""  -1.99%  """
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)

# D48: "1.989" -> "1.993"
$ws.Range("D48").Formula = "=""1.993"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)

# E48: "  -3.87%  " -> "  -3.65%  "
$ws.Range("E48").Formula = "=""  -3.65%  """
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)

# D49: "3.423" -> "3.421"
$ws.Range("D49").Formula = "=""3.421"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)

# E49: "  +0.12%  " -> "  +0.15%  "
$ws.Range("E49").Formula = "=""  +0.15%  """
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)

# D50: "0.06894" -> "0.06900"
$ws.Range("D50").Formula = "=""0.06900"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)

# E50: "  +0.83%  " -> "  +0.96%  "
$ws.Range("E50").Formula = "=""  +0.96%  """
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)

# D51: "110.32" -> "110.46"
$ws.Range("D51").Formula = "=""110.46"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

# E51: "  -0.68%  " -> "  -0.44%  "
$ws.Range("E51").Formula = "=""  -0.44%  """
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)

$excel.CutCopyMode = 0

